# Add repo to repos list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear every hyperlink on the sheet so we can rebuild them in the right
# order (this host's Hyperlinks.Delete() clears the whole sheet's
# collection, not just the calling range's, so everything is re-added
# below).
$ws.Range("A1").Hyperlinks.Delete()

# Capitalize Maria's name (B4)
$ws.Range("B4").Value = "Maria"

# Match the look/format of the other data rows (xlPasteFormats = -4122),
# then fill in row 5's values: the new repo, contributed by Ana.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122) | Out-Null

$ws.Range("A5").Value = "https://github.com/danielhantunes/jiraflow-sample2"
$ws.Range("B5").Value = "Ana"
$ws.Range("C5").Value = "email@example.com"

# Fix Daniel's email (C2): was a personal gmail address, now the generic placeholder
$ws.Range("C2").Value = "email@example.com"

# Re-create all the hyperlinks, in the same order Excel would produce them:
# the untouched ones first, then the brand-new row-5 links, then the
# recreated C2 link (its target text changed so Excel treats it as a new
# hyperlink object appended at the end).
$ws.Range("A2").Hyperlinks.Add($ws.Range("A2"), "https://github.com/danielhantunes/JiraFlow") | Out-Null
$ws.Range("A3").Hyperlinks.Add($ws.Range("A3"), "https://github.com/danielhantunes/jiraflow-sample1") | Out-Null
$ws.Range("C3").Hyperlinks.Add($ws.Range("C3"), "mailto:email@example.com") | Out-Null
$ws.Range("A4").Hyperlinks.Add($ws.Range("A4"), "https://github.com/repoaleatorio/repoaleatorio") | Out-Null
$ws.Range("C4").Hyperlinks.Add($ws.Range("C4"), "mailto:email@example.com") | Out-Null
$ws.Range("A5").Hyperlinks.Add($ws.Range("A5"), "https://github.com/danielhantunes/jiraflow-sample2") | Out-Null
$ws.Range("C5").Hyperlinks.Add($ws.Range("C5"), "mailto:email@example.com") | Out-Null
$ws.Range("C2").Hyperlinks.Add($ws.Range("C2"), "mailto:email@example.com") | Out-Null

# Hyperlinks.Add stamps each cell with a freshly derived (but visually
# identical) style; put them back on the named "Hiperlink" cell style so
# the styles stay the same as the other hyperlinked cells.
$ws.Range("A2").Style = "Hiperlink"
$ws.Range("C2").Style = "Hiperlink"
$ws.Range("A3").Style = "Hiperlink"
$ws.Range("C3").Style = "Hiperlink"
$ws.Range("A4").Style = "Hiperlink"
$ws.Range("C4").Style = "Hiperlink"
$ws.Range("A5").Style = "Hiperlink"
$ws.Range("C5").Style = "Hiperlink"

$ws.Range("B6").Select() | Out-Null
